# Kazakhstan Premier League workbook update (30-05-2024 12:21)
# The source data rows for a handful of matches were re-ordered upstream;
# this swaps the match-detail columns (B:AD) between the affected row
# pairs while leaving column A (the running index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose B:AD contents need to be swapped.
$pairs = @(
    @(99, 100),
    @(107, 108),
    @(143, 144),
    @(150, 151)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $val1 = $cell1.Value()
        $val2 = $cell2.Value()

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}
